# Add columns I (I0) and J (IF) to the sheet, matching the style/format
# of the existing header row and data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
# Copy formatting (bold font, border, centered alignment) from the
# existing last header cell (H1) onto the new header cells, then set
# their text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# --- Data rows (rows 2-44) ---
$iValues = @(4,7,4,5,7,9,9,8,4,6,8,8,7,3,7,8,6,7,7,8,5,9,4,6,6,1,8,4,8,9,7,7,8,6,3,8,6,5,9,4,7,4,3)
$jValues = @(6,8,6,6,8,9,9,8,5,7,9,9,7,6,7,8,6,8,8,8,7,9,5,7,8,2,8,7,8,9,9,9,8,8,3,8,6,6,9,5,7,5,4)

for ($idx = 0; $idx -lt $iValues.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$idx]
    $ws.Cells.Item($row, 10).Value = $jValues[$idx]
}
